$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 10 (question 9 "Какой уровень планирования и организации...") answer text:
# remove the leading sentence "Стратегическое планирование развития направления AI на квартал и год. "
$ws.Range("B10").Value = "Организация спринтов (2 недели), распределение задач между разработчиками. Контроль разработки AI-агентов от идеи до продакшена — архитектура, качество кода, тестирование, деплой."

# Let Excel recompute the (wrapped-text) row height for the shortened content.
$ws.Rows.Item(10).EntireRow.AutoFit()
$ws.Rows.Item(10).RowHeight = 58

# Update the active selection / scrolled view to match the edited cell.
$ws.Range("B10").Select()
$excel.ActiveWindow.ScrollRow = 7
